# The paragraph contains one long bold run:
#   " elementum quis dictum ac, porta ac ante. Fusce tempus ac mauris id
#     cursus. Phasellus a erat nulla. Mauris dolor orci, malesuada auctor
#     dignissim non, posuere nec odio. Etiam hendrerit justo nec diam
#     ullamcorper, nec blandit elit sodales."
#
# We need "Mauris dolor orci" inside that run to also become italic
# (bold+italic), which forces Word to split the single <w:r> into three
# runs: bold-only text, bold+italic "Mauris dolor orci", bold-only text.

$d = $word.ActiveDocument

$target = $d.Content
$target.Find.Execute("Mauris dolor orci", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)

$target.Italic = 1

$d.Save()
